# Auto-generated PowerShell/COM-interop script implementing the MIGS.eu.built.4.0 header-row refactor
# - Inserts organism-group fields (strain/isolate/cultivar/ecotype) starting at column G
# - Shifts existing environmental fields right by one column
# - Appends many new fields (biotic_relationship .. trophic_level) through column BK
# - Re-applies the GREEN/BLUE/YELLOW highlight styles and per-field tooltip comments

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture format (fill color) sources from cells whose formatting is already correct
$blueSrc = $ws.Range("H15")   # existing BLUE "Organism group" cell
$greenSrc = $ws.Range("L15")  # existing GREEN mandatory-field cell
$yellowSrc = $ws.Range("C15") # existing YELLOW optional-field cell

# Re-color G15 (now "strain", part of the BLUE organism group)
$blueSrc.Copy()
$ws.Range("G15").PasteSpecial(-4122)

# Re-color K15 (now "abs_air_humidity", back to GREEN mandatory)
$greenSrc.Copy()
$ws.Range("K15").PasteSpecial(-4122)

# Re-color AL15:BK15 (the new fields) YELLOW (optional)
$ws.Range("AL15:BK15").Select() | Out-Null
$yellowSrc.Copy()
$ws.Range("AL15:BK15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Write the field names for columns G through BK in their new order
$ws.Range("G15").Value2 = 'strain'
$ws.Range("H15").Value2 = 'isolate'
$ws.Range("I15").Value2 = 'cultivar'
$ws.Range("J15").Value2 = 'ecotype'
$ws.Range("K15").Value2 = 'abs_air_humidity'
$ws.Range("L15").Value2 = 'air_temp'
$ws.Range("M15").Value2 = 'build_occup_type'
$ws.Range("N15").Value2 = 'building_setting'
$ws.Range("O15").Value2 = 'carb_dioxide'
$ws.Range("P15").Value2 = 'collection_date'
$ws.Range("Q15").Value2 = 'env_biome'
$ws.Range("R15").Value2 = 'env_feature'
$ws.Range("S15").Value2 = 'env_material'
$ws.Range("T15").Value2 = 'estimated_size'
$ws.Range("U15").Value2 = 'filter_type'
$ws.Range("V15").Value2 = 'geo_loc_name'
$ws.Range("W15").Value2 = 'heat_cool_type'
$ws.Range("X15").Value2 = 'indoor_space'
$ws.Range("Y15").Value2 = 'isol_growth_condt'
$ws.Range("Z15").Value2 = 'lat_lon'
$ws.Range("AA15").Value2 = 'light_type'
$ws.Range("AB15").Value2 = 'num_replicons'
$ws.Range("AC15").Value2 = 'occup_samp'
$ws.Range("AD15").Value2 = 'occupant_dens_samp'
$ws.Range("AE15").Value2 = 'organism_count'
$ws.Range("AF15").Value2 = 'ploidy'
$ws.Range("AG15").Value2 = 'propagation'
$ws.Range("AH15").Value2 = 'rel_air_humidity'
$ws.Range("AI15").Value2 = 'space_typ_state'
$ws.Range("AJ15").Value2 = 'typ_occupant_dens'
$ws.Range("AK15").Value2 = 'ventilation_type'
$ws.Range("AL15").Value2 = 'biotic_relationship'
$ws.Range("AM15").Value2 = 'dew_point'
$ws.Range("AN15").Value2 = 'extrachrom_elements'
$ws.Range("AO15").Value2 = 'health_state'
$ws.Range("AP15").Value2 = 'host'
$ws.Range("AQ15").Value2 = 'host_taxid'
$ws.Range("AR15").Value2 = 'indoor_surf'
$ws.Range("AS15").Value2 = 'isolation_source'
$ws.Range("AT15").Value2 = 'locus_tag_prefix'
$ws.Range("AU15").Value2 = 'pathogenicity'
$ws.Range("AV15").Value2 = 'ref_biomaterial'
$ws.Range("AW15").Value2 = 'samp_collect_device'
$ws.Range("AX15").Value2 = 'samp_mat_process'
$ws.Range("AY15").Value2 = 'samp_size'
$ws.Range("AZ15").Value2 = 'samp_sort_meth'
$ws.Range("BA15").Value2 = 'samp_vol_we_dna_ext'
$ws.Range("BB15").Value2 = 'source_material_id'
$ws.Range("BC15").Value2 = 'subspecf_gen_lin'
$ws.Range("BD15").Value2 = 'substructure_type'
$ws.Range("BE15").Value2 = 'surf_air_cont'
$ws.Range("BF15").Value2 = 'surf_humidity'
$ws.Range("BG15").Value2 = 'surf_material'
$ws.Range("BH15").Value2 = 'surf_moisture'
$ws.Range("BI15").Value2 = 'surf_moisture_ph'
$ws.Range("BJ15").Value2 = 'surf_temp'
$ws.Range("BK15").Value2 = 'trophic_level'

# Drop the old tooltip comments (their text no longer matches the shifted fields)
for ($col = 7; $col -le 38; $col++) {
    $cell = $ws.Cells.Item(15, $col)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete()
    }
}

# Add the correct tooltip comment for every field from G to BK
$ws.Range("G15").AddComment('Organism group

microbial or eukaryotic strain name')
$ws.Range("H15").AddComment('Organism group

Identification or description of the specific individual from which this sample was obtained')
$ws.Range("I15").AddComment('Organism group

Cultivar name - cultivated variety of plant')
$ws.Range("J15").AddComment('Organism group

a population within a given species displaying genetically based, phenotypic traits that reflect adaptation to a local habitat, e.g., Columbia')
$ws.Range("K15").AddComment('actual mass of water vapor - mh20 - present in the air water vapor mixture')
$ws.Range("L15").AddComment('temperature of the air at the time of sampling')
$ws.Range("M15").AddComment('primary function for which a building or discrete part of a building is intended to be used')
$ws.Range("N15").AddComment('location (geography) where a building is set')
$ws.Range("O15").AddComment('carbon dioxide (gas) amount or concentration at the time of sampling')
$ws.Range("P15").AddComment('Time of sampling (single instance or interval, eg., 2008-01-23T19:23:10, 2008-01-23, 2008-01, 2008, 1952-10-21T11:43Z/1952-10-21T17:43Z, 1952-10-21/1953-02-15, 1952-10/1953-02, 1952/1953)')
$ws.Range("Q15").AddComment('Descriptor of the broad ecological context of a sample. Examples include: desert, taiga or deciduous woodland. FAQ, http://trace.ddbj.nig.ac.jp/biosample/faq_e.html#biome-feature-material EnvO (v 2013-06-14) terms can be found via the link: http://www.environmentontology.org/Browse-EnvO')
$ws.Range("R15").AddComment('Descriptor of the local environment. Examples include: harbor, cliff, or lake. FAQ, http://trace.ddbj.nig.ac.jp/biosample/faq_e.html#biome-feature-material EnvO (v 2013-06-14) terms can be found via the link: http://www.environmentontology.org/Browse-EnvO')
$ws.Range("S15").AddComment('Material that was displaced by the sample, or material in which a sample was embedded, prior to the sampling event. Examples include: air, soil, or water. FAQ, http://trace.ddbj.nig.ac.jp/biosample/faq_e.html#biome-feature-material EnvO (v 2013-06-14) terms can be found via the link: http://www.environmentontology.org/Browse-EnvO')
$ws.Range("T15").AddComment('Estimated size of genome')
$ws.Range("U15").AddComment('device which removes solid particulates or airborne molecular contaminants')
$ws.Range("V15").AddComment('Geographical origin of the sample; use the appropriate name from the list, http://www.ddbj.nig.ac.jp/sub/country-e.html. Use a colon to separate the country or ocean from more detailed information about the location, eg "Japan:Kanagawa, Hakone, Lake Ashi" ')
$ws.Range("W15").AddComment('methods of conditioning or heating a room or building')
$ws.Range("X15").AddComment('a distinguishable space within a structure, the purpose for which discrete areas of a building is used')
$ws.Range("Y15").AddComment('Publication reference in the form of pubmed ID, DOI or URL for isolation and growth condition specifications of the organism/material')
$ws.Range("Z15").AddComment('The geographical coordinates of the location where the sample was collected. Specify as decimal degrees latitude and longitude in format "d[d.dddd] N|S d[dd.dddd] W|E", eg, 47.94 N 28.12 W')
$ws.Range("AA15").AddComment('application of light to achieve some practical or aesthetic effect. Lighting includes the use of both artificial light sources such as lamps and light fixtures, as well as natural illumination by capturing daylight. Can also include absence of light')
$ws.Range("AB15").AddComment('Reports the number of replicons in a nuclear genome of eukaryotes, in the genome of a bacterium or archaea or the number of segments in a segmented virus. Always applied to the haploid chromosome count of a eukaryote')
$ws.Range("AC15").AddComment('number of occupants present at time of sample within the given space')
$ws.Range("AD15").AddComment('average number of occupants at time of sampling per square footage')
$ws.Range("AE15").AddComment('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
$ws.Range("AF15").AddComment('The ploidy level of the genome (e.g. allopolyploid, haploid, diploid, triploid, tetraploid). ')
$ws.Range("AG15").AddComment('This field is specific to different taxa. For phage: lytic/lysogenic/temperate/obligately lytic;  for plasmid: incompatibility group;  for eukaryote: asexual/sexual')
$ws.Range("AH15").AddComment('partial vapor and air pressure, density of the vapor and air, or by the actual mass of the vapor and air')
$ws.Range("AI15").AddComment('customary or normal state of the space')
$ws.Range("AJ15").AddComment('customary or normal density of occupants')
$ws.Range("AK15").AddComment('ventilation system used in the sampled premises')
$ws.Range("AL15").AddComment('Free-living or from host (define relationship)')
$ws.Range("AM15").AddComment('temperature to which a given parcel of humid air must be cooled, at constant barometric pressure, for water vapor to condense into water.')
$ws.Range("AN15").AddComment('Plasmids that have significance phenotypic consequence')
$ws.Range("AO15").AddComment('Health or disease status of sample at time of collection')
$ws.Range("AP15").AddComment('The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".')
$ws.Range("AQ15").AddComment('NCBI taxonomy ID of the host, e.g. 9606')
$ws.Range("AR15").AddComment('type of indoor surface')
$ws.Range("AS15").AddComment('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$ws.Range("AT15").AddComment('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html')
$ws.Range("AU15").AddComment('To what is the entity pathogenic')
$ws.Range("AV15").AddComment('Primary publication or genome report in the form of pubmed ID, DOI or URL')
$ws.Range("AW15").AddComment('Method or device employed for collecting sample')
$ws.Range("AX15").AddComment('Processing applied to the sample during or after isolation')
$ws.Range("AY15").AddComment('Amount or size of sample (volume, mass or area) that was collected')
$ws.Range("AZ15").AddComment('method by which samples are sorted')
$ws.Range("BA15").AddComment('volume (mL) or weight (g) of sample processed for DNA extraction')
$ws.Range("BB15").AddComment('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$ws.Range("BC15").AddComment('Information about the genetic distinctness of the lineage (eg., biovar, serovar)')
$ws.Range("BD15").AddComment('substructure or under building is that largely hidden section of the building which is built off the foundations to the ground floor level')
$ws.Range("BE15").AddComment('contaminant identified on surface')
$ws.Range("BF15").AddComment('surfaces: water activity as a function of air and material moisture')
$ws.Range("BG15").AddComment('surface materials at the point of sampling')
$ws.Range("BH15").AddComment('water held on a surface')
$ws.Range("BI15").AddComment('pH measurement of surface')
$ws.Range("BJ15").AddComment('temperature of the surface at the time of sampling')
$ws.Range("BK15").AddComment('Feeding position in food chain (eg., chemolithotroph)')

